# Applies the per-row currentAveragePrice / Leve profit-column refresh
# produced by the scheduled market-data runner (see commit message).
# Each row is addressed as Sheet!H<row>:N<row> matching the
# currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ columns.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 50001224
$ws.Range("I137").Value = 62501092
$ws.Range("J137").Value = 1745.75
$ws.Range("K137").Value = 187503276
$ws.Range("L137").Value = 5237.25
$ws.Range("M137").Value = -187500726
$ws.Range("N137").Value = -10337.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 22203
$ws.Range("I21").Value = 30338.334
$ws.Range("J21").Value = 10000
$ws.Range("K21").Value = 30338.334
$ws.Range("L21").Value = 10000
$ws.Range("M21").Value = -29964.334
$ws.Range("N21").Value = -10748

$ws.Range("H32").Value = 25698.217
$ws.Range("I32").Value = 3026.2856
$ws.Range("J32").Value = 263753.5
$ws.Range("K32").Value = 3026.2856
$ws.Range("L32").Value = 263753.5
$ws.Range("M32").Value = -2739.2856
$ws.Range("N32").Value = -264327.5

$ws.Range("H61").Value = 2926.276
$ws.Range("I61").Value = 2336.6365
$ws.Range("J61").Value = 4779.4287
$ws.Range("K61").Value = 2336.6365
$ws.Range("L61").Value = 4779.4287
$ws.Range("M61").Value = -2124.6365
$ws.Range("N61").Value = -5203.4287

$ws.Range("M86").ClearContents()
$ws.Range("H86").Value = 50000
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 50000
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 50000
$ws.Range("N86").Value = -52372

$ws.Range("M89").ClearContents()
$ws.Range("H89").Value = 50000
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 50000
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 150000
$ws.Range("N89").Value = -161856

$ws.Range("H114").Value = 30838.8
$ws.Range("J114").Value = 30838.8
$ws.Range("L114").Value = 30838.8
$ws.Range("N114").Value = -39516.8

$ws.Range("H136").Value = 2926.276
$ws.Range("I136").Value = 2336.6365
$ws.Range("J136").Value = 4779.4287
$ws.Range("K136").Value = 7009.9095
$ws.Range("L136").Value = 14338.2861
$ws.Range("M136").Value = -4459.9095
$ws.Range("N136").Value = -19438.2861

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1475.125
$ws.Range("I20").Value = 1332.2
$ws.Range("J20").Value = 1713.3334
$ws.Range("K20").Value = 1332.2
$ws.Range("L20").Value = 1713.3334
$ws.Range("M20").Value = -1085.2
$ws.Range("N20").Value = -2207.3334

$ws.Range("H134").Value = 3317.5945
$ws.Range("I134").Value = 2160.739
$ws.Range("K134").Value = 6482.217000000001
$ws.Range("M134").Value = -3947.217000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H80").Value = 28114
$ws.Range("J80").Value = 28114
$ws.Range("L80").Value = 28114
$ws.Range("N80").Value = -30360

$ws.Range("H83").Value = 28114
$ws.Range("J83").Value = 28114
$ws.Range("L83").Value = 84342
$ws.Range("N83").Value = -95574

$ws.Range("H134").Value = 2901.8
$ws.Range("I134").Value = 1551.3636
$ws.Range("J134").Value = 6615.5
$ws.Range("K134").Value = 4654.0908
$ws.Range("L134").Value = 19846.5
$ws.Range("M134").Value = -2119.0908
$ws.Range("N134").Value = -24916.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1116.1842
$ws.Range("I5").Value = 511.21738
$ws.Range("J5").Value = 2043.8
$ws.Range("K5").Value = 1533.65214
$ws.Range("L5").Value = 6131.4
$ws.Range("M5").Value = -1421.65214
$ws.Range("N5").Value = -6355.4

$ws.Range("H16").Value = 760.4
$ws.Range("I16").Value = 700.5
$ws.Range("K16").Value = 2101.5
$ws.Range("M16").Value = -1928.5

$ws.Range("H23").Value = 211.88235
$ws.Range("I23").Value = 80
$ws.Range("J23").Value = 220.125
$ws.Range("K23").Value = 240
$ws.Range("L23").Value = 660.375
$ws.Range("M23").Value = -5
$ws.Range("N23").Value = -1130.375

$ws.Range("H60").Value = 2192.5
$ws.Range("I60").Value = 375
$ws.Range("J60").Value = 2646.875
$ws.Range("K60").Value = 1125
$ws.Range("L60").Value = 7940.625
$ws.Range("M60").Value = -874
$ws.Range("N60").Value = -8442.625

$ws.Range("H109").Value = 1598.3572
$ws.Range("I109").Value = 1364.75
$ws.Range("J109").Value = 3000
$ws.Range("K109").Value = 4094.25
$ws.Range("L109").Value = 9000
$ws.Range("M109").Value = -3054.25
$ws.Range("N109").Value = -11080

$ws.Range("N121").ClearContents()
$ws.Range("H121").Value = 296.66666
$ws.Range("I121").Value = 296.66666
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 889.9999799999999
$ws.Range("L121").Value = 0
$ws.Range("M121").Value = 420.0000200000001

$ws.Range("H135").Value = 1116.1842
$ws.Range("I135").Value = 511.21738
$ws.Range("J135").Value = 2043.8
$ws.Range("K135").Value = 4600.95642
$ws.Range("L135").Value = 18394.2
$ws.Range("M135").Value = -2065.95642
$ws.Range("N135").Value = -23464.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 16521.2
$ws.Range("J103").Value = 16521.2
$ws.Range("L103").Value = 16521.2
$ws.Range("N103").Value = -18865.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("N4").ClearContents()
$ws.Range("H4").Value = 17966.666
$ws.Range("I4").Value = 17966.666
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 17966.666
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -17853.666

$ws.Range("N28").ClearContents()
$ws.Range("H28").Value = 17966.666
$ws.Range("I28").Value = 17966.666
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 17966.666
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -17734.666

$ws.Range("H34").Value = 22416.666
$ws.Range("I34").Value = 34833.332
$ws.Range("J34").Value = 10000
$ws.Range("K34").Value = 34833.332
$ws.Range("L34").Value = 10000
$ws.Range("M34").Value = -34661.332
$ws.Range("N34").Value = -10344

$ws.Range("N37").ClearContents()
$ws.Range("H37").Value = 17966.666
$ws.Range("I37").Value = 17966.666
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 17966.666
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -17859.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 29875
$ws.Range("I28").Value = 50000
$ws.Range("J28").Value = 23166.666
$ws.Range("K28").Value = 50000
$ws.Range("L28").Value = 23166.666
$ws.Range("M28").Value = -49652
$ws.Range("N28").Value = -23862.666

$ws.Range("H31").Value = 7250
$ws.Range("J31").Value = 7250
$ws.Range("L31").Value = 7250
$ws.Range("N31").Value = -7946

$ws.Range("N96").ClearContents()
$ws.Range("H96").Value = 932.6667
$ws.Range("I96").Value = 932.6667
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 932.6667
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = 440.3333

$ws.Range("H132").Value = 2326.9512
$ws.Range("I132").Value = 1991.7407
$ws.Range("J132").Value = 2973.4285
$ws.Range("K132").Value = 5975.2221
$ws.Range("L132").Value = 8920.2855
$ws.Range("M132").Value = -3445.2221
$ws.Range("N132").Value = -13980.2855
